# Schema change: introduce a new "RepresentationVariant" sheet carrying the
# representation-variant fields that used to live inline on ResourceInfo
# (rdf_url/rdf_type/schema_url/schema_type), collapsing them into a single
# "representation_variants" link column, and push the existing "Container"
# sheet out to make room right after it.

$wb = $excel.ActiveWorkbook

# --- 1. ResourceInfo: collapse rdf_url/rdf_type/schema_url/schema_type (D:G)
#        into a single "representation_variants" column (D). ---------------
$resourceInfo = $wb.Worksheets.Item("ResourceInfo")
$resourceInfo.Range("D1").Value = "representation_variants"
$resourceInfo.Range("E1:G1").ClearContents()

# --- 2. Turn the existing "Container" sheet into "RepresentationVariant" --
#        (keeps its sheetId / position, matching the diff) and give it the
#        new header row. ----------------------------------------------------
$repVariant = $wb.Worksheets.Item("Container")
$repVariant.Name = "RepresentationVariant"
$repVariant.Range("A1").Value = "url"
$repVariant.Range("B1").Value = "media_type"
$repVariant.Range("C1").Value = "encoding_format"
$repVariant.Range("D1").Value = "size"

# --- 3. Add a fresh "Container" sheet right after it, restoring the old --
#        single-column Container layout. ------------------------------------
$container = $wb.Worksheets.Add($null, $repVariant)
$container.Name = "Container"
$container.Range("A1").Value = "contains_pids"
